$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Style = "Normal"
$ws.Range("A1").Value = "a"
$ws.Range("A1").Font.Color = 255
$ws.Range("A1").Interior.Color = 65535
